# Update loading_percent values for Case_2_251 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 15.71019891993584
    "C2" = 12.60612085043198
    "E2" = 17.08198617088914
    "F2" = 35.22975175446599
    "G2" = 26.41881194302372
    "H2" = 13.67201851351358
    "J2" = 7.24095835845102
    "L2" = 12.87326636577617
    "N2" = 17.28768089724858
    "O2" = 20.50361336558248
    "B3" = 15.22634282696266
    "C3" = 12.56307518449164
    "E3" = 17.1113394368775
    "F3" = 35.23918306540212
    "G3" = 26.42529109241418
    "H3" = 13.71473794688386
    "J3" = 7.235799127588362
    "L3" = 12.84191017979012
    "N3" = 17.3261427164476
    "O3" = 20.5622702292523
    "B4" = 14.92288244016674
    "C4" = 12.53699818222509
    "E4" = 17.13187412044775
    "F4" = 35.25370401714996
    "G4" = 26.4388276588916
    "H4" = 13.74332395013261
    "J4" = 7.232694890052205
    "L4" = 12.82432990772062
    "N4" = 17.35151842529197
    "O4" = 20.60313747116825
    "B5" = 14.79780130592586
    "C5" = 12.52646516292385
    "E5" = 17.14087394816019
    "F5" = 35.2618166063728
    "G5" = 26.44674228369536
    "H5" = 13.75556504344364
    "J5" = 7.23144610562865
    "L5" = 12.81759120510813
    "N5" = 17.36230249619365
    "O5" = 20.62100809660506
    "B6" = 14.77695163452194
    "C6" = 12.52472195174953
    "E6" = 17.14240652502939
    "F6" = 35.26329625284674
    "G6" = 26.44820116479796
    "H6" = 13.75763341513389
    "J6" = 7.231239732292885
    "L6" = 12.81649807662219
    "N6" = 17.36411997569092
    "O6" = 20.62404890309713
    "B7" = 14.92120104335626
    "C7" = 12.53685574516385
    "E7" = 17.13199293700943
    "F7" = 35.2538045391776
    "G7" = 26.4389246957892
    "H7" = 13.74348664108761
    "J7" = 7.232677982534967
    "L7" = 12.82423729864798
    "N7" = 17.35166206735554
    "O7" = 20.60337355778
    "B8" = 15.54479329973182
    "C8" = 12.59120826723882
    "E8" = 17.09158602086469
    "F8" = 35.23119206321494
    "G8" = 26.41905964712711
    "H8" = 13.68625893321592
    "J8" = 7.239166407569855
    "L8" = 12.86211059859497
    "N8" = 17.30057766251384
    "O8" = 20.52282947464056
    "B9" = 16.70976550838757
    "C9" = 12.70038599166667
    "E9" = 17.03226670233251
    "F9" = 35.25607074764595
    "G9" = 26.45607904505107
    "H9" = 13.59274745024003
    "J9" = 7.252384467586276
    "L9" = 12.94942696543947
    "N9" = 17.21433539402393
    "O9" = 20.40351002551464
    "B10" = 17.52155524986746
    "C10" = 12.78189647515431
    "E10" = 17.00081125793063
    "F10" = 35.31640913608582
    "G10" = 26.52963316975224
    "H10" = 13.53547488161007
    "J10" = 7.262384367438582
    "L10" = 13.02121469309788
    "N10" = 17.15942649704317
    "O10" = 20.33956221177141
    "B11" = 17.87969543203666
    "C11" = 12.81919794873366
    "E11" = 16.98912955362399
    "F11" = 35.3529342562137
    "G11" = 26.57312160091742
    "H11" = 13.51190677926979
    "J11" = 7.266993496238914
    "L11" = 13.05545734762038
    "N11" = 17.13627406032647
    "O11" = 20.31565175644176
    "B12" = 18.01359671421504
    "C12" = 12.83334959802712
    "E12" = 16.98508326269931
    "F12" = 35.36806342509588
    "G12" = 26.59102453283775
    "H12" = 13.50333994487828
    "J12" = 7.268747239631327
    "L12" = 13.0686455188487
    "N12" = 17.12776875948543
    "O12" = 20.30734447694353
    "B13" = 17.98483682354465
    "C13" = 12.83030070187382
    "E13" = 16.98593793179023
    "F13" = 35.36474749520785
    "G13" = 26.58710514293497
    "H13" = 13.50516904098724
    "J13" = 7.268369172015579
    "L13" = 13.06579548319511
    "N13" = 17.12958888301226
    "O13" = 20.3091003381985
    "B14" = 17.8907466333477
    "C14" = 12.82036171633898
    "E14" = 16.98878910412991
    "F14" = 35.35415297722638
    "G14" = 26.57456578004818
    "H14" = 13.5111948052051
    "J14" = 7.267137611706446
    "L14" = 13.05653795297264
    "N14" = 17.13556907577927
    "O14" = 20.31495332681346
    "B15" = 17.83288663845753
    "C15" = 12.81427707280289
    "E15" = 16.99058465039616
    "F15" = 35.34783231207152
    "G15" = 26.5670716574408
    "H15" = 13.51493238225248
    "J15" = 7.266384326629516
    "L15" = 13.05089605012888
    "N15" = 17.13926622411593
    "O15" = 20.31863580682628
    "B16" = 17.49791546095506
    "C16" = 12.77946273109444
    "E16" = 17.00162751512332
    "F16" = 35.31420418741644
    "G16" = 26.52699232704849
    "H16" = 13.53706516870286
    "J16" = 7.262084336101854
    "L16" = 13.01900819089713
    "N16" = 17.16097625872271
    "O16" = 20.34122924437668
    "B17" = 17.28948011273143
    "C17" = 12.75815819076538
    "E17" = 17.00907458815508
    "F17" = 35.29589444411144
    "G17" = 26.50496902704614
    "H17" = 13.5512798590574
    "J17" = 7.259461633079843
    "L17" = 12.99984745834687
    "N17" = 17.17476193133687
    "O17" = 20.35641792750939
    "B18" = 17.16855048833392
    "C18" = 12.74592560712837
    "E18" = 17.01360530025587
    "F18" = 35.28621805132877
    "G18" = 26.49324630055708
    "H18" = 13.55968969002148
    "J18" = 7.257958799843735
    "L18" = 12.98897642618104
    "N18" = 17.18286297886434
    "O18" = 20.36564147879565
    "B19" = 17.12743037063062
    "C19" = 12.74178767734293
    "E19" = 17.01518181896261
    "F19" = 35.28308881982489
    "G19" = 26.48943958699578
    "H19" = 13.56257727461176
    "J19" = 7.257450949660996
    "L19" = 12.98532159952513
    "N19" = 17.18563539377469
    "O19" = 20.3688480569879
    "B20" = 17.31177724893045
    "C20" = 12.76042393866628
    "E20" = 17.00825623793159
    "F20" = 35.2977551187194
    "G20" = 26.50721573850257
    "H20" = 13.54974246787078
    "J20" = 7.259740238839724
    "L20" = 13.00187170496036
    "N20" = 17.17327663657216
    "O20" = 20.35475060266785
    "B21" = 17.91843070627251
    "C21" = 12.82328036533819
    "E21" = 16.98794141000506
    "F21" = 35.35722968479422
    "G21" = 26.57821002397663
    "H21" = 13.50941517527008
    "J21" = 7.267499126012026
    "L21" = 13.05925116601971
    "N21" = 17.13380544115144
    "O21" = 20.31321386975011
    "B22" = 18.30485129300484
    "C22" = 12.86451190547604
    "E22" = 16.97686348373546
    "F22" = 35.40366051453216
    "G22" = 26.63296704272528
    "H22" = 13.48514525085863
    "J22" = 7.27261863990461
    "L22" = 13.09803804901291
    "N22" = 17.10953577188119
    "O22" = 20.29042268099974
    "B23" = 18.09956710452879
    "C23" = 12.8424938557092
    "E23" = 16.98257496848214
    "F23" = 35.37819048106432
    "G23" = 26.60298041631891
    "H23" = 13.49790751555592
    "J23" = 7.269881902610097
    "L23" = 13.07722141226707
    "N23" = 17.12234940123243
    "O23" = 20.30218756895984
    "B24" = 17.30170012277061
    "C24" = 12.75939954478075
    "E24" = 17.0086254375195
    "F24" = 35.29691125912722
    "G24" = 26.50619707555237
    "H24" = 13.55043678229545
    "J24" = 7.259614265597908
    "L24" = 13.00095609174041
    "N24" = 17.17394759167833
    "O24" = 20.35550286979175
    "B25" = 16.40177341572056
    "C25" = 12.67060187297364
    "E25" = 17.04618311658615
    "F25" = 35.24193880756838
    "G25" = 26.4379176356096
    "H25" = 13.61603910398016
    "J25" = 7.248757013194948
    "L25" = 12.92444036844791
    "N25" = 17.23617879178139
    "O25" = 20.43163552828682
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
